# Table 16 SI formatting pass:
#  - header row: shrink italic header font 11pt -> 10pt, restyle 3 header
#    cells (blank the first, add a "2017" second line under Norganizations,
#    rename Norganizations_norm -> Norganizations' with a "(normalized)"
#    second line)
#  - data rows: right-align the country-code column, normalize numeric
#    cells to two decimal places

$d = $word.ActiveDocument
$t = $d.Tables(1)

# ---- Header row (row 1) -------------------------------------------------

# Cell 1: "COUNTRIES" -> blank, font 22 half-pts -> 20 half-pts (11pt->10pt)
$cell = $t.Cell(1, 1)
$cell.Range.Font.Size = 10
$content = $d.Range($cell.Range.Start, $cell.Range.End - 1)
$content.Text = ""

# Cell 2: "Norganizations" stays, font shrinks, gains a line-break + "2017"
$cell = $t.Cell(1, 2)
$content = $d.Range($cell.Range.Start, $cell.Range.End - 1)
$len = $content.End - $content.Start
$content.Text = "Norganizations" + [char]11 + "2017"
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$headPart = $d.Range($cellStart, $cellStart + $len)
$headPart.Font.Size = 10
$tailPart = $d.Range($cellStart + $len, $cellEnd - 1)
$tailPart.Font.Size = 10

# Cell 3: "Norganizations_norm" -> "Norganizations'" + line-break + "(normalized)"
$cell = $t.Cell(1, 3)
$content = $d.Range($cell.Range.Start, $cell.Range.End - 1)
$newHead = "Norganizations'"
$content.Text = $newHead + [char]11 + "(normalized)"
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$headPart = $d.Range($cellStart, $cellStart + $newHead.Length)
$headPart.Font.Size = 10
$tailPart = $d.Range($cellStart + $newHead.Length, $cellEnd - 1)
$tailPart.Font.Size = 10

# Cell 4: "CO.MANAGEMENT" text unchanged, just the font shrinks
$cell = $t.Cell(1, 4)
$content = $d.Range($cell.Range.Start, $cell.Range.End - 1)
$content.Font.Size = 10

# ---- Data rows ------------------------------------------------------------

# Right-align the country-code (column 1) paragraph on every data row.
for ($r = 2; $r -le $t.Rows.Count; $r++) {
    $t.Cell($r, 1).Range.ParagraphFormat.Alignment = 2
}

# Normalize numeric cell text to two decimal places (values that were
# integers or had only one decimal place gain the trailing zero(es)).
$numericEdits = @{
    2  = @{ 2 = "1.00";  3 = "0.00";  4 = "0.00" };   # BE
    3  = @{ 2 = "2.00" };                              # DK
    4  = @{ 2 = "14.00" };                             # DE
    5  = @{ 2 = "7.00" };                              # EE
    6  = @{ 2 = "5.00";  3 = "0.10";  4 = "0.10" };   # IE
    7  = @{ 2 = "40.00"; 3 = "1.00";  4 = "1.00" };   # ES
    8  = @{ 2 = "21.00" };                             # FR
    9  = @{ 2 = "4.00" };                              # LV
    10 = @{ 2 = "3.00" };                              # LT
    11 = @{ 2 = "12.00" };                             # NL
    12 = @{ 2 = "11.00" };                             # PL
    13 = @{ 2 = "15.00" };                             # PT
    15 = @{ 2 = "6.00" };                              # SE
}

foreach ($row in $numericEdits.Keys) {
    $cols = $numericEdits[$row]
    foreach ($col in $cols.Keys) {
        $cell = $t.Cell($row, $col)
        $content = $d.Range($cell.Range.Start, $cell.Range.End - 1)
        $content.Text = $cols[$col]
    }
}
